$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.681.64'
$ws.Range("E2").Value = '  -2.13%  '
$ws.Range("D3").Value = '2.352.97'
$ws.Range("E3").Value = '  -0.51%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.638'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.52%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.622'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.02'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.66%  '
$ws.Range("E11").Value = '  -1.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.43'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.72%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.996'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.01%  '
$ws.Range("E14").Value = '  +0.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.11'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.92%  '
$ws.Range("D16").Value = '2.705.38'
$ws.Range("E16").Value = '  -0.55%  '
$ws.Range("D17").Value = '2.348.26'
$ws.Range("E17").Value = '  -4.57%  '
$ws.Range("D18").Value = '42.634.91'
$ws.Range("E18").Value = '  -2.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.43%  '
$ws.Range("E20").Value = '  -2.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '76.43'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.69'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '263.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.31'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -9.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.99'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.75%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.42'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.37%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.70'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.39%  '
$ws.Range("E29").Value = '  -1.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '175.18'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.18%  '
$ws.Range("E31").Value = '  -3.64%  '
$ws.Range("E32").Value = '  -3.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '35.23'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -9.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.04'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.79%  '
$ws.Range("E35").Value = '  -0.49%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.55'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -8.18%  '
$ws.Range("E37").Value = '  +4.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0356'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.89%  '
$ws.Range("E39").Value = '  -8.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.81'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.238'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.24%  '
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '69.87'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '121.62'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.88%  '
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '92.50'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +22.66%  '
$ws.Range("E47").Value = '  -7.60%  '
$ws.Range("E48").Value = '  -2.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.16'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.47%  '
$ws.Range("E50").Value = '  -4.12%  '
$ws.Range("E51").Value = '  -0.28%  '
